# Corrección a Diebold Mariano y revisión de Cap1
# Updates the P_valores and Estadisticos_DM matrices with corrected values.

$wb = $excel.ActiveWorkbook

# --- P_valores sheet ---
$ws2 = $wb.Worksheets.Item("P_valores")
$ws2.Range("C2").Value = 0.6413660283590135
$ws2.Range("D2").Value = 0.9547186884386549
$ws2.Range("E2").Value = 0.9118127404151573
$ws2.Range("F2").Value = 0.1289416893904598
$ws2.Range("B3").Value = 0.6413660283590135
$ws2.Range("D3").Value = 0.4353446046275207
$ws2.Range("E3").Value = 0.2949445627359824
$ws2.Range("F3").Value = 0.4351323334701944
$ws2.Range("B4").Value = 0.9547186884386549
$ws2.Range("C4").Value = 0.4353446046275207
$ws2.Range("E4").Value = 0.9714403183110845
$ws2.Range("F4").Value = 0.2470400207514738
$ws2.Range("B5").Value = 0.9118127404151573
$ws2.Range("C5").Value = 0.2949445627359824
$ws2.Range("D5").Value = 0.9714403183110845
$ws2.Range("F5").Value = 0.1835358555734827
$ws2.Range("B6").Value = 0.1289416893904598
$ws2.Range("C6").Value = 0.4351323334701944
$ws2.Range("D6").Value = 0.2470400207514738
$ws2.Range("E6").Value = 0.1835358555734827

# --- Estadisticos_DM sheet ---
$ws3 = $wb.Worksheets.Item("Estadisticos_DM")
$ws3.Range("C2").Value = -0.472297189299867
$ws3.Range("D2").Value = 0.05743305881543417
$ws3.Range("E2").Value = 0.1120335421137764
$ws3.Range("F2").Value = -1.577546623584738
$ws3.Range("B3").Value = 0.472297189299867
$ws3.Range("D3").Value = 0.7945828633713284
$ws3.Range("E3").Value = 1.072899315494204
$ws3.Range("F3").Value = -0.7949554796708516
$ws3.Range("B4").Value = -0.05743305881543417
$ws3.Range("C4").Value = -0.7945828633713284
$ws3.Range("E4").Value = 0.03621145594854015
$ws3.Range("F4").Value = -1.18919573139237
$ws3.Range("B5").Value = -0.1120335421137764
$ws3.Range("C5").Value = -1.072899315494204
$ws3.Range("D5").Value = -0.03621145594854015
$ws3.Range("F5").Value = -1.37316458277484
$ws3.Range("B6").Value = 1.577546623584738
$ws3.Range("C6").Value = 0.7949554796708516
$ws3.Range("D6").Value = 1.18919573139237
$ws3.Range("E6").Value = 1.37316458277484
